# Insert 6 new data rows (new market-day readings) immediately above the
# existing row 337 block, pushing the old rows 337-354 down to 343-360.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 blank rows starting at row 337.
$ws.Range("A337:A342").EntireRow.Insert()

# Constant values shared by every row in this sub-block.
$mercadoId = 3
$mercado   = "Femacal de La Calera"
$region    = "Coquimbo"
$codreg    = 5
$catId     = 100112027
$categoria = "Melón"
$unidad    = "`$/unidad"
$kgUnidad  = 1
$clasif    = "Hortaliza"
$fecha     = "2022-01-24"
$origen    = "Provincia de Talca"

# Row 337: Calameño - Extra
$r = 337
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $catId
$ws.Cells.Item($r, 7).Value = $categoria
$ws.Cells.Item($r, 8).Value = "Calameño"
$ws.Cells.Item($r, 9).Value = "Extra"
$ws.Cells.Item($r, 10).Value = 300
$ws.Cells.Item($r, 11).Value = 1000
$ws.Cells.Item($r, 12).Value = 1000
$ws.Cells.Item($r, 13).Value = 1000
$ws.Cells.Item($r, 14).Value = $unidad
$ws.Cells.Item($r, 15).Value = $origen
$ws.Cells.Item($r, 16).Value = 1000
$ws.Cells.Item($r, 17).Value = $kgUnidad
$ws.Cells.Item($r, 18).Value = $clasif

# Row 338: Calameño - Primera
$r = 338
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $catId
$ws.Cells.Item($r, 7).Value = $categoria
$ws.Cells.Item($r, 8).Value = "Calameño"
$ws.Cells.Item($r, 9).Value = "Primera"
$ws.Cells.Item($r, 10).Value = 300
$ws.Cells.Item($r, 11).Value = 700
$ws.Cells.Item($r, 12).Value = 700
$ws.Cells.Item($r, 13).Value = 700
$ws.Cells.Item($r, 14).Value = $unidad
$ws.Cells.Item($r, 15).Value = $origen
$ws.Cells.Item($r, 16).Value = 700
$ws.Cells.Item($r, 17).Value = $kgUnidad
$ws.Cells.Item($r, 18).Value = $clasif

# Row 339: Calameño - Segunda
$r = 339
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $catId
$ws.Cells.Item($r, 7).Value = $categoria
$ws.Cells.Item($r, 8).Value = "Calameño"
$ws.Cells.Item($r, 9).Value = "Segunda"
$ws.Cells.Item($r, 10).Value = 350
$ws.Cells.Item($r, 11).Value = 500
$ws.Cells.Item($r, 12).Value = 500
$ws.Cells.Item($r, 13).Value = 500
$ws.Cells.Item($r, 14).Value = $unidad
$ws.Cells.Item($r, 15).Value = $origen
$ws.Cells.Item($r, 16).Value = 500
$ws.Cells.Item($r, 17).Value = $kgUnidad
$ws.Cells.Item($r, 18).Value = $clasif

# Row 340: Tuna - Extra
$r = 340
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $catId
$ws.Cells.Item($r, 7).Value = $categoria
$ws.Cells.Item($r, 8).Value = "Tuna"
$ws.Cells.Item($r, 9).Value = "Extra"
$ws.Cells.Item($r, 10).Value = 380
$ws.Cells.Item($r, 11).Value = 1000
$ws.Cells.Item($r, 12).Value = 1000
$ws.Cells.Item($r, 13).Value = 1000
$ws.Cells.Item($r, 14).Value = $unidad
$ws.Cells.Item($r, 15).Value = $origen
$ws.Cells.Item($r, 16).Value = 1000
$ws.Cells.Item($r, 17).Value = $kgUnidad
$ws.Cells.Item($r, 18).Value = $clasif

# Row 341: Tuna - Primera
$r = 341
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $catId
$ws.Cells.Item($r, 7).Value = $categoria
$ws.Cells.Item($r, 8).Value = "Tuna"
$ws.Cells.Item($r, 9).Value = "Primera"
$ws.Cells.Item($r, 10).Value = 300
$ws.Cells.Item($r, 11).Value = 700
$ws.Cells.Item($r, 12).Value = 700
$ws.Cells.Item($r, 13).Value = 700
$ws.Cells.Item($r, 14).Value = $unidad
$ws.Cells.Item($r, 15).Value = $origen
$ws.Cells.Item($r, 16).Value = 700
$ws.Cells.Item($r, 17).Value = $kgUnidad
$ws.Cells.Item($r, 18).Value = $clasif

# Row 342: Tuna - Segunda
$r = 342
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $catId
$ws.Cells.Item($r, 7).Value = $categoria
$ws.Cells.Item($r, 8).Value = "Tuna"
$ws.Cells.Item($r, 9).Value = "Segunda"
$ws.Cells.Item($r, 10).Value = 380
$ws.Cells.Item($r, 11).Value = 500
$ws.Cells.Item($r, 12).Value = 500
$ws.Cells.Item($r, 13).Value = 500
$ws.Cells.Item($r, 14).Value = $unidad
$ws.Cells.Item($r, 15).Value = $origen
$ws.Cells.Item($r, 16).Value = 500
$ws.Cells.Item($r, 17).Value = $kgUnidad
$ws.Cells.Item($r, 18).Value = $clasif
